# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to the H:N (price/profit) columns
# across multiple worksheets, per the scheduled-runner price refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 307.14285  # H2: 91.666664 -> 307.14285
$ws.Cells.Item(2, 10).Value = 1600  # J2: 0 -> 1600
$ws.Cells.Item(2, 12).Value = 1600  # L2: 0 -> 1600
$ws.Cells.Item(2, 14).Value = -1826  # N2: (new) -> -1826
$ws.Cells.Item(6, 8).Value = 35  # H6: 16529.166 -> 35
$ws.Cells.Item(6, 9).Value = 40  # I6: 19832 -> 40
$ws.Cells.Item(6, 11).Value = 120  # K6: 59496 -> 120
$ws.Cells.Item(6, 13).Value = -8  # M6: -59384 -> -8
$ws.Cells.Item(38, 8).Value = 1662.2  # H38: 3240.8 -> 1662.2
$ws.Cells.Item(38, 9).Value = 9.5  # I38: 10.75 -> 9.5
$ws.Cells.Item(38, 10).Value = 4141.25  # J38: 5394.1665 -> 4141.25
$ws.Cells.Item(38, 11).Value = 28.5  # K38: 32.25 -> 28.5
$ws.Cells.Item(38, 12).Value = 12423.75  # L38: 16182.4995 -> 12423.75
$ws.Cells.Item(38, 13).Value = 343.5  # M38: 339.75 -> 343.5
$ws.Cells.Item(38, 14).Value = -13167.75  # N38: -16926.4995 -> -13167.75
$ws.Cells.Item(41, 8).Value = 408.8889  # H41: 430.83334 -> 408.8889
$ws.Cells.Item(41, 9).Value = 326.85715  # I41: 298.25 -> 326.85715
$ws.Cells.Item(41, 11).Value = 326.85715  # K41: 298.25 -> 326.85715
$ws.Cells.Item(41, 13).Value = 113.14285  # M41: 141.75 -> 113.14285
$ws.Cells.Item(43, 8).Value = 482.33334  # H43: 500 -> 482.33334
$ws.Cells.Item(43, 9).Value = 482.33334  # I43: 500 -> 482.33334
$ws.Cells.Item(43, 11).Value = 482.33334  # K43: 500 -> 482.33334
$ws.Cells.Item(43, 13).Value = -413.33334  # M43: -431 -> -413.33334
$ws.Cells.Item(101, 8).Value = 1000  # H101: 427.66666 -> 1000
$ws.Cells.Item(101, 9).Value = 1000  # I101: 549 -> 1000
$ws.Cells.Item(101, 10).Value = 0  # J101: 185 -> 0
$ws.Cells.Item(101, 11).Value = 3000  # K101: 1647 -> 3000
$ws.Cells.Item(101, 12).Value = 0  # L101: 555 -> 0
$ws.Cells.Item(101, 13).ClearContents()  # M101: -25 -> (removed)
$ws.Cells.Item(101, 14).ClearContents()  # N101: -3799 -> (removed)

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(21, 8).Value = 0  # H21: 550 -> 0
$ws.Cells.Item(21, 9).Value = 0  # I21: 600 -> 0
$ws.Cells.Item(21, 10).Value = 0  # J21: 500 -> 0
$ws.Cells.Item(21, 11).Value = 0  # K21: 600 -> 0
$ws.Cells.Item(21, 12).ClearContents()  # L21: 500 -> (removed)
$ws.Cells.Item(21, 13).ClearContents()  # M21: -226 -> (removed)
$ws.Cells.Item(21, 14).ClearContents()  # N21: -1248 -> (removed)
$ws.Cells.Item(63, 8).Value = 6864.875  # H63: 7559.857 -> 6864.875
$ws.Cells.Item(63, 10).Value = 6709.857  # J63: 7494.8335 -> 6709.857
$ws.Cells.Item(63, 12).Value = 6709.857  # L63: 7494.8335 -> 6709.857
$ws.Cells.Item(63, 14).Value = -8081.857  # N63: -8866.833500000001 -> -8081.857
$ws.Cells.Item(66, 8).Value = 6864.875  # H66: 7559.857 -> 6864.875
$ws.Cells.Item(66, 10).Value = 6709.857  # J66: 7494.8335 -> 6709.857
$ws.Cells.Item(66, 12).Value = 33549.285  # L66: 37474.1675 -> 33549.285
$ws.Cells.Item(66, 14).Value = -40413.285  # N66: -44338.1675 -> -40413.285
$ws.Cells.Item(92, 8).Value = 50996.668  # H92: 50122.5 -> 50996.668
$ws.Cells.Item(92, 10).Value = 50996.668  # J92: 50122.5 -> 50996.668
$ws.Cells.Item(92, 12).Value = 50996.668  # L92: 50122.5 -> 50996.668
$ws.Cells.Item(92, 14).Value = -55988.668  # N92: -55114.5 -> -55988.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(12, 8).Value = 105  # H12: 425.75 -> 105
$ws.Cells.Item(12, 9).Value = 105  # I12: 236.33333 -> 105
$ws.Cells.Item(12, 10).Value = 0  # J12: 994 -> 0
$ws.Cells.Item(12, 11).Value = 105  # K12: 236.33333 -> 105
$ws.Cells.Item(12, 12).Value = 0  # L12: 994 -> 0
$ws.Cells.Item(12, 13).ClearContents()  # M12: -68.33332999999999 -> (removed)
$ws.Cells.Item(12, 14).ClearContents()  # N12: -1330 -> (removed)
$ws.Cells.Item(82, 8).Value = 34569.8  # H82: 21166.666 -> 34569.8
$ws.Cells.Item(82, 10).Value = 50283  # J82: 41500 -> 50283
$ws.Cells.Item(82, 12).Value = 50283  # L82: 41500 -> 50283
$ws.Cells.Item(82, 14).Value = -51049  # N82: -42266 -> -51049
$ws.Cells.Item(85, 8).Value = 34569.8  # H85: 21166.666 -> 34569.8
$ws.Cells.Item(85, 10).Value = 50283  # J85: 41500 -> 50283
$ws.Cells.Item(85, 12).Value = 50283  # L85: 41500 -> 50283
$ws.Cells.Item(85, 14).Value = -52935  # N85: -44152 -> -52935
$ws.Cells.Item(99, 8).Value = 1773.8  # H99: 1677.3334 -> 1773.8
$ws.Cells.Item(99, 9).Value = 1773.8  # I99: 1677.3334 -> 1773.8
$ws.Cells.Item(99, 11).Value = 1773.8  # K99: 1677.3334 -> 1773.8
$ws.Cells.Item(99, 13).Value = -275.8  # M99: -179.3334 -> -275.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(6, 8).Value = 150505  # H6: 14000000 -> 150505
$ws.Cells.Item(6, 9).Value = 1010  # I6: 14000000 -> 1010
$ws.Cells.Item(6, 10).Value = 300000  # J6: 0 -> 300000
$ws.Cells.Item(6, 11).Value = 1010  # K6: 14000000 -> 1010
$ws.Cells.Item(6, 12).Value = 300000  # L6: 0 -> 300000
$ws.Cells.Item(6, 13).Value = -897  # M6: -13999887 -> -897
$ws.Cells.Item(6, 14).Value = -300226  # N6: (new) -> -300226
$ws.Cells.Item(7, 8).Value = 374.55554  # H7: 408 -> 374.55554
$ws.Cells.Item(7, 9).Value = 358.875  # I7: 392.66666 -> 358.875
$ws.Cells.Item(7, 11).Value = 358.875  # K7: 392.66666 -> 358.875
$ws.Cells.Item(7, 13).Value = -245.875  # M7: -279.66666 -> -245.875
$ws.Cells.Item(17, 8).Value = 0  # H17: 20 -> 0
$ws.Cells.Item(17, 9).Value = 0  # I17: 20 -> 0
$ws.Cells.Item(17, 11).Value = 0  # K17: 20 -> 0
$ws.Cells.Item(17, 13).ClearContents()  # M17: 154 -> (removed)
$ws.Cells.Item(25, 8).Value = 0  # H25: 1000 -> 0
$ws.Cells.Item(25, 9).Value = 0  # I25: 1000 -> 0
$ws.Cells.Item(25, 11).Value = 0  # K25: 1000 -> 0
$ws.Cells.Item(25, 13).ClearContents()  # M25: -826 -> (removed)
$ws.Cells.Item(31, 8).Value = 7534.154  # H31: 6956.2 -> 7534.154
$ws.Cells.Item(31, 9).Value = 4680.75  # I31: 4427.222 -> 4680.75
$ws.Cells.Item(31, 10).Value = 12099.6  # J31: 10749.667 -> 12099.6
$ws.Cells.Item(31, 11).Value = 4680.75  # K31: 4427.222 -> 4680.75
$ws.Cells.Item(31, 12).Value = 12099.6  # L31: 10749.667 -> 12099.6
$ws.Cells.Item(31, 13).Value = -4385.75  # M31: -4132.222 -> -4385.75
$ws.Cells.Item(31, 14).Value = -12689.6  # N31: -11339.667 -> -12689.6
$ws.Cells.Item(34, 8).Value = 7534.154  # H34: 6956.2 -> 7534.154
$ws.Cells.Item(34, 9).Value = 4680.75  # I34: 4427.222 -> 4680.75
$ws.Cells.Item(34, 10).Value = 12099.6  # J34: 10749.667 -> 12099.6
$ws.Cells.Item(34, 11).Value = 4680.75  # K34: 4427.222 -> 4680.75
$ws.Cells.Item(34, 12).Value = 12099.6  # L34: 10749.667 -> 12099.6
$ws.Cells.Item(34, 13).Value = -4478.75  # M34: -4225.222 -> -4478.75
$ws.Cells.Item(34, 14).Value = -12503.6  # N34: -11153.667 -> -12503.6
$ws.Cells.Item(93, 8).Value = 53069  # H93: 41326.168 -> 53069
$ws.Cells.Item(93, 9).Value = 54603.5  # I93: 39319 -> 54603.5
$ws.Cells.Item(93, 10).Value = 50000  # J93: 43333.332 -> 50000
$ws.Cells.Item(93, 11).Value = 54603.5  # K93: 39319 -> 54603.5
$ws.Cells.Item(93, 12).Value = 50000  # L93: 43333.332 -> 50000
$ws.Cells.Item(93, 13).Value = -52731.5  # M93: -37447 -> -52731.5
$ws.Cells.Item(93, 14).Value = -53744  # N93: -47077.332 -> -53744

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1478.2222  # H5: 1615.1428 -> 1478.2222
$ws.Cells.Item(5, 9).Value = 1103  # I5: 1155.5 -> 1103
$ws.Cells.Item(5, 10).Value = 1665.8334  # J5: 1799 -> 1665.8334
$ws.Cells.Item(5, 11).Value = 3309  # K5: 3466.5 -> 3309
$ws.Cells.Item(5, 12).Value = 4997.5002  # L5: 5397 -> 4997.5002
$ws.Cells.Item(5, 13).Value = -3197  # M5: -3354.5 -> -3197
$ws.Cells.Item(5, 14).Value = -5221.5002  # N5: -5621 -> -5221.5002
$ws.Cells.Item(61, 8).Value = 833  # H61: 500 -> 833
$ws.Cells.Item(61, 10).Value = 999.5  # J61: 0 -> 999.5
$ws.Cells.Item(61, 12).Value = 2998.5  # L61: 0 -> 2998.5
$ws.Cells.Item(61, 14).Value = -3428.5  # N61: (new) -> -3428.5
$ws.Cells.Item(135, 8).Value = 1478.2222  # H135: 1615.1428 -> 1478.2222
$ws.Cells.Item(135, 9).Value = 1103  # I135: 1155.5 -> 1103
$ws.Cells.Item(135, 10).Value = 1665.8334  # J135: 1799 -> 1665.8334
$ws.Cells.Item(135, 11).Value = 9927  # K135: 10399.5 -> 9927
$ws.Cells.Item(135, 12).Value = 14992.5006  # L135: 16191 -> 14992.5006
$ws.Cells.Item(135, 13).Value = -7392  # M135: -7864.5 -> -7392
$ws.Cells.Item(135, 14).Value = -20062.5006  # N135: -21261 -> -20062.5006

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(20, 8).Value = 0  # H20: 26666.666 -> 0
$ws.Cells.Item(20, 10).Value = 0  # J20: 26666.666 -> 0
$ws.Cells.Item(20, 12).ClearContents()  # L20: 26666.666 -> (removed)
$ws.Cells.Item(20, 14).ClearContents()  # N20: -27156.666 -> (removed)
$ws.Cells.Item(98, 8).Value = 18019.5  # H98: 23678.666 -> 18019.5
$ws.Cells.Item(98, 10).Value = 18019.5  # J98: 23678.666 -> 18019.5
$ws.Cells.Item(98, 12).Value = 18019.5  # L98: 23678.666 -> 18019.5
$ws.Cells.Item(98, 14).Value = -24009.5  # N98: -29668.666 -> -24009.5
$ws.Cells.Item(132, 8).Value = 7714.857  # H132: 8077.5386 -> 7714.857
$ws.Cells.Item(132, 9).Value = 6923.6924  # I132: 7250.6665 -> 6923.6924
$ws.Cells.Item(132, 11).Value = 20771.0772  # K132: 21751.9995 -> 20771.0772
$ws.Cells.Item(132, 13).Value = -18241.0772  # M132: -19221.9995 -> -18241.0772

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 150  # H55: 934 -> 150
$ws.Cells.Item(55, 9).Value = 150  # I55: 1225 -> 150
$ws.Cells.Item(55, 10).Value = 0  # J55: 497.5 -> 0
$ws.Cells.Item(55, 11).Value = 150  # K55: 1225 -> 150
$ws.Cells.Item(55, 12).Value = 0  # L55: 497.5 -> 0
$ws.Cells.Item(55, 13).ClearContents()  # M55: -1052 -> (removed)
$ws.Cells.Item(55, 14).ClearContents()  # N55: -843.5 -> (removed)
$ws.Cells.Item(101, 8).Value = 21111  # H101: 22222 -> 21111
$ws.Cells.Item(101, 10).Value = 21111  # J101: 22222 -> 21111
$ws.Cells.Item(101, 12).Value = 21111  # L101: 22222 -> 21111
$ws.Cells.Item(101, 14).Value = -27601  # N101: -28712 -> -27601
$ws.Cells.Item(104, 8).Value = 6035  # H104: 7370 -> 6035
$ws.Cells.Item(104, 10).Value = 6035  # J104: 7370 -> 6035
$ws.Cells.Item(104, 12).Value = 6035  # L104: 7370 -> 6035
$ws.Cells.Item(104, 14).Value = -13023  # N104: -14358 -> -13023
$ws.Cells.Item(134, 8).Value = 30000.5  # H134: 0 -> 30000.5
$ws.Cells.Item(134, 10).Value = 30000.5  # J134: 0 -> 30000.5
$ws.Cells.Item(134, 12).Value = 30000.5  # L134: 0 -> 30000.5
$ws.Cells.Item(134, 14).Value = -40140.5  # N134: (new) -> -40140.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(3, 8).Value = 0  # H3: 7500 -> 0
$ws.Cells.Item(3, 9).Value = 0  # I3: 7500 -> 0
$ws.Cells.Item(3, 11).Value = 0  # K3: 7500 -> 0
$ws.Cells.Item(3, 13).ClearContents()  # M3: -7386 -> (removed)
$ws.Cells.Item(105, 8).Value = 29866.5  # H105: 23900 -> 29866.5
$ws.Cells.Item(105, 10).Value = 29866.5  # J105: 23900 -> 29866.5
$ws.Cells.Item(105, 12).Value = 29866.5  # L105: 23900 -> 29866.5
$ws.Cells.Item(105, 14).Value = -36854.5  # N105: -30888 -> -36854.5
$ws.Cells.Item(132, 8).Value = 6624.875  # H132: 7235 -> 6624.875
$ws.Cells.Item(132, 9).Value = 3833.3333  # I132: 6125.7144 -> 3833.3333
$ws.Cells.Item(132, 10).Value = 14999.5  # J132: 15000 -> 14999.5
$ws.Cells.Item(132, 11).Value = 11499.9999  # K132: 18377.1432 -> 11499.9999
$ws.Cells.Item(132, 12).Value = 44998.5  # L132: 45000 -> 44998.5
$ws.Cells.Item(132, 13).Value = -8969.999899999999  # M132: -15847.1432 -> -8969.999899999999
$ws.Cells.Item(132, 14).Value = -50058.5  # N132: -50060 -> -50058.5
